$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.065.94"
$ws.Range("E2").Value = "  -1.17%  "
$ws.Range("D3").Value = "1.808.87"
$ws.Range("E3").Value = "  -2.18%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.611"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  +0.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.38"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.59%  "
$ws.Range("E9").Value = "  +6.06%  "
$ws.Range("E10").Value = "  -0.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1000"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("D12").Value = "2.070.97"
$ws.Range("E12").Value = "  -2.19%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.816.55"
$ws.Range("E13").Value = "  -1.69%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.663"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.58%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "11.06"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.61%  "
$ws.Range("D17").Value = "35.047.10"
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "237.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.38%  "
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("E24").Value = "  +2.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.21%  "
$ws.Range("E29").Value = "  +18.42%  "
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("E31").Value = "  +4.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0558"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.37%  "
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.695"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.62%  "
$ws.Range("E36").Value = "  +5.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "92.07"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.76%  "
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("D39").Value = "1.313.66"
$ws.Range("E39").Value = "  -1.91%  "
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.80%  "
$ws.Range("E43").Value = "  -2.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.07%  "
$ws.Range("E45").Value = "  -2.26%  "
$ws.Range("E46").Value = "  +4.48%  "
$ws.Range("E47").Value = "  -1.82%  "
$ws.Range("D48").Value = "1.988.30"
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("E49").Value = "  +0.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0664"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "99.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.68%  "
